$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header B1: "locacalizacion" -> "localizacion"
$ws.Range("B1").Value = "localizacion"

# Remove the value from E2 (previously 1)
$ws.Range("E2").ClearContents()

# Update selection to active cell E2
$ws.Range("E2").Select()
